$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 38 (rows 38..56 shift down to 40..58).
$ws.Range("A38:A39").EntireRow.Insert()

# Copy the date column's number format so the new D38/D39 cells keep the
# same date formatting as the rest of the column (they now sit above the
# old row 38, which moved down to row 40).
$ws.Range("D38:D39").NumberFormat = $ws.Range("D40").NumberFormat

# New row 38: Melón, Extra quality, Región de O'Higgins.
$ws.Cells.Item(38, 1).Value = 8
$ws.Cells.Item(38, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(38, 3).Value = "Coquimbo"
$ws.Cells.Item(38, 4).Value = 44566
$ws.Cells.Item(38, 5).Value = 4
$ws.Cells.Item(38, 6).Value = 100112027
$ws.Cells.Item(38, 7).Value = "Melón"
$ws.Cells.Item(38, 8).Value = "Tuna"
$ws.Cells.Item(38, 9).Value = "Extra"
$ws.Cells.Item(38, 10).Value = 6000
$ws.Cells.Item(38, 11).Value = 1100
$ws.Cells.Item(38, 12).Value = 1200
$ws.Cells.Item(38, 13).Value = 1150
$ws.Cells.Item(38, 14).Value = "$/unidad"
$ws.Cells.Item(38, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(38, 16).Value = 1150
$ws.Cells.Item(38, 17).Value = 1
$ws.Cells.Item(38, 18).Value = "Hortaliza"

# New row 39: Melón, Primera quality, Región de O'Higgins.
$ws.Cells.Item(39, 1).Value = 8
$ws.Cells.Item(39, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(39, 3).Value = "Coquimbo"
$ws.Cells.Item(39, 4).Value = 44566
$ws.Cells.Item(39, 5).Value = 4
$ws.Cells.Item(39, 6).Value = 100112027
$ws.Cells.Item(39, 7).Value = "Melón"
$ws.Cells.Item(39, 8).Value = "Tuna"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 5000
$ws.Cells.Item(39, 11).Value = 950
$ws.Cells.Item(39, 12).Value = 1000
$ws.Cells.Item(39, 13).Value = 975
$ws.Cells.Item(39, 14).Value = "$/unidad"
$ws.Cells.Item(39, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(39, 16).Value = 975
$ws.Cells.Item(39, 17).Value = 1
$ws.Cells.Item(39, 18).Value = "Hortaliza"
